# Refresh crypto price/volume columns (D, E) to the latest scrape.
# D-column cells that look like plain decimal numbers (single "." , all digits)
# are forced back to Text (NumberFormat "@") then restyled to "Normal" so the
# written value keeps its original text type/style (matches source t="inlineStr"),
# instead of Excel auto-coercing "508.40" -> the number 508.4.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "56.264.52"
$ws.Range("E2").Value = "  -1.73%  "

$ws.Range("D3").Value = "3.001.82"
$ws.Range("E3").Value = "  +0.18%  "

$ws.Range("E4").Value = "  +0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "508.40"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.30%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "138.49"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.10%  "

$ws.Range("E7").Value = "  -0.04%  "

$ws.Range("E8").Value = "  -0.14%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "7.11"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.39%  "

$ws.Range("E10").Value = "  -0.30%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.368"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.21%  "

$ws.Range("D12").Value = "3.520.30"
$ws.Range("E12").Value = "  +0.39%  "

$ws.Range("E13").Value = "  -0.44%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "25.41"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.36%  "

$ws.Range("E15").Value = "  +1.39%  "

$ws.Range("D16").Value = "56.261.56"
$ws.Range("E16").Value = "  -1.76%  "

$ws.Range("D17").Value = "3.005.41"
$ws.Range("E17").Value = "  +0.47%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "5.94"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.21%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.93"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.58%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.01"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.98%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "332.82"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +4.04%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.999"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.07%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.497"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.08%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "64.80"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.51%  "

$ws.Range("D25").Value = "3.132.67"
$ws.Range("E25").Value = "  +0.56%  "

$ws.Range("E26").Value = "  +1.57%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.999"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.15%  "

$ws.Range("D28").Value = "0.0₃0935"
$ws.Range("E28").Value = "  +4.99%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.34"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.17%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.86"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.03%  "

$ws.Range("E31").Value = "  +0.66%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "20.34"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.02%  "

$ws.Range("E33").Value = "  -0.03%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "152.75"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.38%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.43"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.40%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "26.50"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +8.62%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.81"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.72%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.22"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.70%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0660"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.18%  "

$ws.Range("D40").Value = "3.043.44"
$ws.Range("E40").Value = "  +0.57%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "36.36"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.98%  "

$ws.Range("E42").Value = "  +0.02%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.78"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.29%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.656"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.53%  "

$ws.Range("D45").Value = "2.200.13"
$ws.Range("E45").Value = "  +0.88%  "

$ws.Range("E46").Value = "  -2.66%  "

$ws.Range("E47").Value = "  +2.46%  "

# Rows 48-49: coin order swapped (Cosmos now listed before ONDO); new price/volume too.
$ws.Range("B48").Value = "Cosmos"
$ws.Range("C48").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "5.83"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.85%  "

$ws.Range("B49").Value = "ONDO"
$ws.Range("C49").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.922"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.76%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "19.48"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.60%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0849"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.12%  "

